# Update "想去人数" (interest count) figures for two entries that appear
# on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $ws.Range("F2").Value = 423
        $ws.Range("F3").Value = 2898
        $ws.Range("F5").Value = 53
    }
    elseif ($name -eq "全部类型") {
        $ws.Range("F2").Value = 423
        $ws.Range("F7").Value = 2898
        $ws.Range("F10").Value = 53
    }
}
